# Daily attendance processing - rotate the "Recorded By" (column G) list
# so the first recorder in the comma-separated list is moved to the end.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2

    if ($v -ne $null -and $v -like "*, *") {
        $parts = $v -split ", "
        if ($parts.Length -gt 1) {
            $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
            $cell.Value2 = $rotated
        }
    }
}
